# Cardshifter card list: add the 'Enhancement' card type formulas to F22:F30
# (replacing the blank/"null" shared-string placeholder that was there before),
# widen column F to fit the new text, and leave the selection on F22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F22 gets its own (non-shared) formula producing the literal string 'Enhancement'
$ws.Range("F22").Formula = '="''Enhancement''"'

# F23:F30 share one formula (mirrors the existing 'Bio'/'B0T' shared-formula pattern
# already used higher up in the same column).
$ws.Range("F23:F30").Formula = '="''Enhancement''"'

# Column F needs to be widened now that it holds "'Enhancement'" instead of blank.
$ws.Columns("F").ColumnWidth = 12.5

# Leave the visible selection on F22, as in the authored change.
$ws.Range("F22").Select() | Out-Null
